$d = $word.ActiveDocument

# Locate the paragraph that ends with "Percent correct (need to add to database)"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Percent correct (need to add to database)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Collapse to a point just before the paragraph mark (End - 1) so the
    # insertion lands strictly inside this paragraph's story, rather than on
    # the shared boundary with the following paragraph (which would cause
    # InsertXML to clobber a neighboring paragraph instead of adding a new one).
    $insertAt = $d.Range($target.Range.End - 1, $target.Range.End - 1)

    $xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:t>Show a different percent for easy, medium, and hard problems</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $insertAt.InsertXML($xml) | Out-Null
}
